$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 LogisticRegression(C=3, max_iter=1000, penalty=''l1'',
                                    random_state=42, solver=''liblinear''))])'
$ws.Range("B2").Value = 0.657142857142857
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': None, ''model__solver'': ''liblinear'', ''model__penalty'': ''l1'', ''model__class_weight'': None, ''model__C'': 3}'
$ws.Range("D2").Value = 0.4285714285714285
$ws.Range("E2").Value = '[1 0 0 1 0 0 1 1 0 1 0 0]'
$ws.Range("F2").Value = '[1 0 1 0 1 1 1 1 1 0 1 1]'
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.6720982142857143
$ws.Range("I2").Value = 0.02757616013048205
$ws.Range("J2").Value = 0.5688988095238096
$ws.Range("K2").Value = 0.0576685283626349

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 LogisticRegression(C=0.0001, max_iter=1000, random_state=42,
                                    solver=''liblinear''))])'
$ws.Range("B3").Value = 0.6285714285714284
$ws.Range("C3").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': None, ''model__solver'': ''liblinear'', ''model__penalty'': ''l2'', ''model__class_weight'': None, ''model__C'': 0.0001}'
$ws.Range("D3").Value = 0.625
$ws.Range("E3").Value = '[1 0 1 0 0 0 0 1 1 0 1 1]'
$ws.Range("F3").Value = '[1 1 1 1 1 0 1 1 0 1 1 1]'
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.6718081435472741
$ws.Range("I3").Value = 0.02626674109451626
$ws.Range("J3").Value = 0.5620427881297446
$ws.Range("K3").Value = 0.05862149506263768

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 LogisticRegression(C=0.001, class_weight=''balanced'',
                                    max_iter=1000, random_state=42,
                                    solver=''liblinear''))])'
$ws.Range("B4").Value = 0.6285714285714284
$ws.Range("C4").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': None, ''model__solver'': ''liblinear'', ''model__penalty'': ''l2'', ''model__class_weight'': ''balanced'', ''model__C'': 0.001}'
$ws.Range("D4").Value = 0.75
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[1 0 1 1 1 1 1 0 1 0 0 1]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.6454648526077097
$ws.Range("I4").Value = 0.03256913678233541
$ws.Range("J4").Value = 0.5413454270597127
$ws.Range("K4").Value = 0.09148660837376747
